# Update "想去人数" (F column) figures and one E-column time range,
# per gh-pages data refresh at commit 456a3b4.
$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Cells.Item(2, 6).Value = 1560
$wsExpo.Cells.Item(3, 6).Value = 931
$wsExpo.Cells.Item(4, 6).Value = 927
$wsExpo.Cells.Item(5, 6).Value = 528
$wsExpo.Cells.Item(6, 6).Value = 7949
$wsExpo.Cells.Item(7, 6).Value = 141
$wsExpo.Cells.Item(8, 6).Value = 38
$wsExpo.Cells.Item(9, 6).Value = 1939
$wsExpo.Cells.Item(10, 6).Value = 5743
$wsExpo.Cells.Item(11, 6).Value = 580
$wsExpo.Cells.Item(12, 6).Value = 240
$wsExpo.Cells.Item(13, 6).Value = 309
$wsExpo.Cells.Item(14, 6).Value = 8111
$wsExpo.Cells.Item(14, 5).Value = "2024.07.20 09:30-07.22 17:00"
$wsExpo.Cells.Item(15, 6).Value = 9441
$wsExpo.Cells.Item(16, 6).Value = 1158
$wsExpo.Cells.Item(17, 6).Value = 937
$wsExpo.Cells.Item(18, 6).Value = 4572
$wsExpo.Cells.Item(19, 6).Value = 705
$wsExpo.Cells.Item(20, 6).Value = 281
$wsExpo.Cells.Item(21, 6).Value = 87
$wsExpo.Cells.Item(22, 6).Value = 294
$wsExpo.Cells.Item(24, 6).Value = 1221
$wsExpo.Cells.Item(26, 6).Value = 1723
$wsExpo.Cells.Item(27, 6).Value = 761
$wsExpo.Cells.Item(28, 6).Value = 991
$wsExpo.Cells.Item(29, 6).Value = 25
$wsExpo.Cells.Item(30, 6).Value = 1916
$wsExpo.Cells.Item(32, 6).Value = 487
$wsExpo.Cells.Item(33, 6).Value = 2378
$wsExpo.Cells.Item(35, 6).Value = 123
$wsExpo.Cells.Item(36, 6).Value = 1515
$wsExpo.Cells.Item(37, 6).Value = 76
$wsExpo.Cells.Item(38, 6).Value = 1332
$wsExpo.Cells.Item(39, 6).Value = 9
$wsExpo.Cells.Item(40, 6).Value = 821
$wsExpo.Cells.Item(41, 6).Value = 529
$wsExpo.Cells.Item(42, 6).Value = 203
$wsExpo.Cells.Item(43, 6).Value = 59
$wsExpo.Cells.Item(45, 6).Value = 521
$wsExpo.Cells.Item(46, 6).Value = 17
$wsExpo.Cells.Item(47, 6).Value = 864
$wsExpo.Cells.Item(49, 6).Value = 4122

$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Cells.Item(2, 6).Value = 5409

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Cells.Item(2, 6).Value = 1560
$wsAll.Cells.Item(3, 6).Value = 931
$wsAll.Cells.Item(4, 6).Value = 927
$wsAll.Cells.Item(5, 6).Value = 528
$wsAll.Cells.Item(6, 6).Value = 141
$wsAll.Cells.Item(7, 6).Value = 38
$wsAll.Cells.Item(11, 6).Value = 5743
$wsAll.Cells.Item(12, 6).Value = 580
$wsAll.Cells.Item(13, 6).Value = 8111
$wsAll.Cells.Item(13, 5).Value = "2024.07.20 09:30-07.22 17:00"
$wsAll.Cells.Item(14, 6).Value = 9441
$wsAll.Cells.Item(16, 6).Value = 1158
$wsAll.Cells.Item(17, 6).Value = 937
$wsAll.Cells.Item(18, 6).Value = 705
$wsAll.Cells.Item(19, 6).Value = 281
$wsAll.Cells.Item(20, 6).Value = 87
$wsAll.Cells.Item(21, 6).Value = 294
$wsAll.Cells.Item(24, 6).Value = 1221
$wsAll.Cells.Item(26, 6).Value = 761
$wsAll.Cells.Item(27, 6).Value = 991
$wsAll.Cells.Item(28, 6).Value = 25
$wsAll.Cells.Item(29, 6).Value = 1916
$wsAll.Cells.Item(31, 6).Value = 487
$wsAll.Cells.Item(32, 6).Value = 2378
$wsAll.Cells.Item(33, 6).Value = 76
$wsAll.Cells.Item(40, 6).Value = 529
$wsAll.Cells.Item(42, 6).Value = 203
$wsAll.Cells.Item(43, 6).Value = 59
$wsAll.Cells.Item(45, 6).Value = 521
$wsAll.Cells.Item(46, 6).Value = 17
$wsAll.Cells.Item(47, 6).Value = 864
$wsAll.Cells.Item(49, 6).Value = 4122

